$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: title cell text is unchanged; its shared-string index shifts
#     automatically once the now-unused "Claude_complex_perline_output2.csv"
#     string is dropped from the pool (happens naturally as row 26 is rewritten).

# --- Row 26: replace the old "Line-by-line / Claude_complex..." entry with
#     the new "Whole-scan" example row. Written in A, D, E, B, C order so the
#     shared-string table lands in the same sequence as the target workbook.
$ws.Range("A26").Value = "Whole-scan"
$ws.Range("D26").Value = "`"Here's the recreated content of the table from the image: Arrêté le vingt cinq janvier 1920 Demarche servais 10² vingt six janvier Deflandre Gustave Désiré Tubize 22 8b 1919 Dubblemanns Delphine 898 510 - 107/1919 - - 10³ d Ypersiel Julien Ghislain Nivelles 15 8b 1918 Elexance Célestin 250 - 250 150/1919 - - Arrêté le vingt six janvier 1920 servais 10⁴ vingt sept janvier Cordeur Louis Braine l'Alleud 21 janvier 1919 Cordeur Antoine 300 - 300 111/1919 - - 27 janvier 1920 Arrêté le vingt sept janvier 1920 servais 10⁵ vingt huit janvier Decock Léonie Tubize 22 février 1919 Meuris David - 586 586 121/1920 - - 10⁶ d Decock Adèle - 10 8b 1915 - - - - 122/1920 - - 11 d Delendries Anastasie Nivelles 12 7b 1919 Delendries Delphine 1590 1395 195 - 11 mars 1920 8 avril 1920 Sommier 1920 169 Arrêté le vingt huit janvier 1920 servais Arrêté le vingt neuf janvier 1920 servais 11² trente janvier Rousseau Charles Gn Nivelles 18 mars 1919 Rousseau Louis 1500 - 1500 123/1919 - - 11³ d Dedoncker Vital Tubize 22 7b 1919 Brassemans Léocadie 4687 1536 3151 124/1919 - -`""
$ws.Range("E26").Value = "`"Arrêté le vingt cinq janvier 1920 Demarche servais 10² vingt six janvier Deflandre Gustave Désiré Tubize 22 8b 1919 Dubblemanns Delphine 898 510 - 107/1919 - - 10³ d Ypersiel Julien Ghislain Nivelles 15 8b 1918 Elexance Célestin 250 - 250 150/1919 - - Arrêté le vingt six janvier 1920 servais 10⁴ vingt sept janvier Cordeur Louis Braine l'Alleud 21 janvier 1919 Cordeur Antoine 300 - 300 111/1919 - - 27 janvier 1920 Arrêté le vingt sept janvier 1920 servais 10⁵ vingt huit janvier Decock Léonie Tubize 22 février 1919 Meuris David - 586 586 121/1920 - - 10⁶ d Decock Adèle - 10 8b 1915 - - - - 122/1920 - - 11 d Delendries Anastasie Nivelles 12 7b 1919 Delendries Delphine 1590 1395 195 - 11 mars 1920 8 avril 1920 Sommier 1920 169 Arrêté le vingt huit janvier 1920 servais Arrêté le vingt neuf janvier 1920 servais 11² trente janvier Rousseau Charles Gn Nivelles 18 mars 1919 Rousseau Louis 1500 - 1500 123/1919 - - 11³ d Dedoncker Vital Tubize 22 7b 1919 Brassemans Léocadie 4687 1536 3151 124/1919 - -`""
$ws.Range("B26").Value = "claude_one_example_whole_output.csv // new_transcription16.txt"
$ws.Range("C26").Value = "Example 17 // ID 16"

# Row 26 grows tall to fit the big recreated-table text.
$ws.Rows.Item(26).RowHeight = 221

# --- Rows 27:34 were blank filler rows ("Line-by-line" only) - remove them.
$ws.Range("A27:E34").EntireRow.Delete()

# --- Rows 14, 15, 18 shrink back from the old 68pt height to 34pt.
$ws.Rows.Item(14).RowHeight = 34
$ws.Rows.Item(15).RowHeight = 34
$ws.Rows.Item(18).RowHeight = 34

# --- View state: scrolled down, zoomed in, selection on D32.
$ws.Application.ActiveWindow.Zoom = 134
$ws.Application.ActiveWindow.ScrollRow = 24
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D32").Select()
